$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.212.88'
$ws.Range('E2').Value = '  -3.69%  '
$ws.Range('D3').Value = '3.039.14'
$ws.Range('E3').Value = '  -3.35%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '''531.64'
$ws.Range('E5').Value = '  -5.85%  '
$ws.Range('D6').Value = '''130.94'
$ws.Range('E6').Value = '  -9.96%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '3.036.21'
$ws.Range('E8').Value = '  -3.27%  '
$ws.Range('E9').Value = '  -1.79%  '
$ws.Range('D10').Value = '''0.150'
$ws.Range('E10').Value = '  -2.96%  '
$ws.Range('E11').Value = '  -10.66%  '
$ws.Range('D12').Value = '''0.448'
$ws.Range('E12').Value = '  -4.34%  '
$ws.Range('D13').Value = '''0.0000221'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('D14').Value = '''33.79'
$ws.Range('E14').Value = '  -8.61%  '
$ws.Range('D15').Value = '3.487.14'
$ws.Range('E15').Value = '  -4.13%  '
$ws.Range('D16').Value = '62.131.20'
$ws.Range('E16').Value = '  -3.72%  '
$ws.Range('E17').Value = '  -2.34%  '
$ws.Range('D18').Value = '3.044.25'
$ws.Range('E18').Value = '  -2.73%  '
$ws.Range('D19').Value = '''6.48'
$ws.Range('E19').Value = '  -5.61%  '
$ws.Range('D20').Value = '''472.91'
$ws.Range('E20').Value = '  -8.33%  '
$ws.Range('D21').Value = '''13.10'
$ws.Range('E21').Value = '  -6.73%  '
$ws.Range('D22').Value = '''0.689'
$ws.Range('E22').Value = '  -4.02%  '
$ws.Range('D23').Value = '''7.04'
$ws.Range('E23').Value = '  -6.05%  '
$ws.Range('D24').Value = '''78.11'
$ws.Range('E24').Value = '  -1.10%  '
$ws.Range('D25').Value = '''11.82'
$ws.Range('E25').Value = '  -8.84%  '
$ws.Range('E26').Value = '  -0.18%  '
$ws.Range('D27').Value = '''2.64'
$ws.Range('E27').Value = '  -6.78%  '
$ws.Range('E28').Value = '  -11.19%  '
$ws.Range('D29').Value = '''1.00'
$ws.Range('E29').Value = '  +0.48%  '
$ws.Range('D30').Value = '''25.42'
$ws.Range('E30').Value = '  -4.51%  '
$ws.Range('D31').Value = '''1.82'
$ws.Range('E31').Value = '  -16.12%  '
$ws.Range('D32').Value = '''1.09'
$ws.Range('E32').Value = '  -4.99%  '
$ws.Range('E33').Value = '  -9.46%  '
$ws.Range('D34').Value = '''56.50'
$ws.Range('E34').Value = '  +4.54%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').Value = '''5.16'
$ws.Range('E35').Value = '  -4.73%  '
$ws.Range('B36').Value = 'Filecoin'
$ws.Range('C36').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D36').Value = '''5.81'
$ws.Range('E36').Value = '  -4.86%  '
$ws.Range('D37').Value = '''465.89'
$ws.Range('E37').Value = '  -15.08%  '
$ws.Range('D38').Value = '3.064.12'
$ws.Range('E38').Value = '  -2.46%  '
$ws.Range('D39').Value = '''0.0385'
$ws.Range('E39').Value = '  -11.58%  '
$ws.Range('D40').Value = '''0.0777'
$ws.Range('E40').Value = '  -6.22%  '
$ws.Range('D41').Value = '''7.92'
$ws.Range('E41').Value = '  -4.45%  '
$ws.Range('E42').Value = '  -9.69%  '
$ws.Range('D43').Value = '''2.52'
$ws.Range('E43').Value = '  -8.20%  '
$ws.Range('E44').Value = '  +0.03%  '
$ws.Range('E45').Value = '  -9.15%  '
$ws.Range('D46').Value = '''1.99'
$ws.Range('E46').Value = '  -10.69%  '
$ws.Range('B47').Value = 'PEPE'
$ws.Range('C47').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D47').Value = '0.0₃0512'
$ws.Range('E47').Value = '  -2.23%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').Value = '''23.88'
$ws.Range('E48').Value = '  -6.57%  '
$ws.Range('E49').Value = '  -2.46%  '
$ws.Range('D50').Value = '''115.38'
$ws.Range('E50').Value = '  -4.63%  '
$ws.Range('D51').Value = '''1.94'
$ws.Range('E51').Value = '  -8.12%  '
